$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill rows 2 through 36 with the bundle-diameter pull data.
# Column A: incrementing pull number (1..35)
# Column B: cable size "3C#9"
# Column C: Local/Express "EXPRESS"
# Column D: From "543+00"
# Column E: To "553+00"
for ($i = 2; $i -le 36; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
    $ws.Cells.Item($i, 2).Value = "3C#9"
    $ws.Cells.Item($i, 3).Value = "EXPRESS"
    $ws.Cells.Item($i, 4).Value = "543+00"
    $ws.Cells.Item($i, 5).Value = "553+00"
}

# Match the final on-screen selection/scroll state left behind by the edit.
$ws.Range("H36").Select()

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A21:A31").Select()

# Re-activate Sheet1 so it is the tab shown/selected when the file is reopened.
$ws.Activate()

Write-Output "done"
